$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Question")
$ws.Activate()

# The order in which brand-new text is written matters: Excel appends newly
# seen strings to the shared-string table in first-use order, so we mirror
# the exact sequence the author typed the replacement quiz content in.

# Row 2: "cat" question -> "DDD stands for?"
$ws.Range("C2").Value = "DDD stands for?"
$ws.Range("D2").Value = "MCQ"
$ws.Range("F2").Value = "ACB"
$ws.Range("G2").Value = "HEF"
$ws.Range("E2").Value = "DDD"
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

# Row 3: "dog" question -> "Who can have a alpha?"
$ws.Range("E3").Value = "Manager"
$ws.Range("F3").Value = "Driver"
$ws.Range("G3").Value = "You"
$ws.Range("H3").Value = "SME"
$ws.Range("C3").Value = "Who can have a alpha?"
$ws.Range("D3").Value = "MCQ"
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0

# Row 4: "tiger" question -> "Which offshore team handle issue related a driver?"
$ws.Range("E4").Value = "Backoffice"
$ws.Range("H4").Value = "Operations"
$ws.Range("C4").Value = "Which offshore team handle issue related a driver?"
$ws.Range("D4").Value = "MCQ"
$ws.Range("F4").Value = "Order"
$ws.Range("G4").Value = "Vehicle"
$ws.Range("I4").ClearContents()
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0
$ws.Range("N4").ClearContents()

# Row 5: "lion" question -> "How are all SMEs for Vehicle?"
$ws.Range("C5").Value = "How are all SMEs for Vehicle?"
$ws.Range("D5").Value = "MCA"
$ws.Range("E5").Value = "Jamie"
$ws.Range("F5").Value = "Anitha"
$ws.Range("G5").Value = "Maria"
$ws.Range("H5").Value = "Vlad"
$ws.Range("I5").Value = "Prasanna"
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 0

# Update the active selection on the Question sheet to L6
$ws.Range("L6").Select()
